$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "From" value of the R30 rule row (C10) from 18 to 1
$ws.Range("C10").Value = 1
